# Updates cryptos list values (price/volume columns) and fixes three rows
# whose order got reshuffled (coin name/link/price/volume moved together),
# per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.584.44"
$ws.Range("D3").Value = "1.827.66"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'317.65"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.5391"
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("D8").Value = "'0.4008"
$ws.Range("E8").Value = "  +6.49%  "
$ws.Range("D9").Value = "'0.07772"
$ws.Range("E9").Value = "  +4.90%  "
$ws.Range("D10").Value = "'1.120"
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("D11").Value = "'41.98"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "'21.34"
$ws.Range("E12").Value = "  +4.07%  "
$ws.Range("D13").Value = "'6.341"
$ws.Range("D14").Value = "'7.644"
$ws.Range("E14").Value = "  +5.94%  "
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "1.826.88"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "'0.00001094"
$ws.Range("E17").Value = "  +3.56%  "
$ws.Range("D18").Value = "'90.01"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "'0.06593"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "'17.77"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "'6.080"
$ws.Range("E22").Value = "  +3.29%  "
$ws.Range("D23").Value = "28.595.19"
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("D24").Value = "'11.20"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "'2.267"
$ws.Range("E25").Value = "  +8.20%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'20.87"
$ws.Range("E26").Value = "  +3.00%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.464"
$ws.Range("E27").Value = "  +8.17%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'158.22"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "2.038.62"
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("D30").Value = "'124.46"
$ws.Range("E30").Value = "  +2.94%  "
$ws.Range("D31").Value = "'1.135"
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("D32").Value = "'0.1119"
$ws.Range("E32").Value = "  +5.52%  "
$ws.Range("D33").Value = "'5.701"
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("D34").Value = "'0.07503"
$ws.Range("E34").Value = "  +16.50%  "
$ws.Range("D35").Value = "'3.652"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").Value = "'0.2255"
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("D37").Value = "'0.02360"
$ws.Range("E37").Value = "  +3.13%  "
$ws.Range("D38").Value = "'8.980"
$ws.Range("E38").Value = "  +6.39%  "
$ws.Range("D39").Value = "'5.219"
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6327"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'11.39"
$ws.Range("E41").Value = "  +2.54%  "
$ws.Range("D42").Value = "'1.191"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "'1.406"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").Value = "'13.56"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").Value = "'0.5918"
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("D47").Value = "'3.707"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").Value = "'125.44"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("E49").Value = "  +4.54%  "
$ws.Range("D50").Value = "'1.199"
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").Value = "'0.06915"
$ws.Range("E51").Value = "  +1.45%  "
